$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Hôtel Duo"
$ws.Range("B4").Value = "US$686"
$ws.Range("C4").Value = "'8.4"
$ws.Range("D4").Value = "Very Good"
$ws.Range("E4").Value = "'850"

$ws.Range("A5").Value = "St Christopher's Inn Paris - Gare du Nord"
$ws.Range("B5").Value = "US$95"
$ws.Range("C5").Value = "'7.2"
$ws.Range("D5").Value = "Good"
$ws.Range("E5").Value = "'9,248"

$ws.Range("A6").Value = "St Christopher's Inn Paris - Canal"
$ws.Range("B6").Value = "US$86"
$ws.Range("C6").Value = "'7.3"
$ws.Range("D6").Value = "Good"
$ws.Range("E6").Value = "'6,519"

$ws.Range("A7").Value = "Generator Paris"
$ws.Range("B7").Value = "US$101"
$ws.Range("C7").Value = "'7.3"
$ws.Range("D7").Value = "Good"
$ws.Range("E7").Value = "'9,523"

$ws.Range("A8").Value = "The People - Paris Nation"
$ws.Range("B8").Value = "US$121"
$ws.Range("C8").Value = "'8.4"
$ws.Range("D8").Value = "Very Good"
$ws.Range("E8").Value = "'6,122"

$ws.Range("A9").Value = "The People - Paris Bercy"
$ws.Range("B9").Value = "US$143"
$ws.Range("C9").Value = "'8.8"
$ws.Range("D9").Value = "Excellent"
$ws.Range("E9").Value = "'5,976"

$ws.Range("A10").Value = "Hotel de L'Empereur by Malone"
$ws.Range("B10").Value = "US$741"
$ws.Range("C10").Value = "'8.6"
$ws.Range("D10").Value = "Excellent"
$ws.Range("E10").Value = "'992"

$ws.Range("A11").Value = "Le Regent Montmartre by Hiphophostels"
$ws.Range("B11").Value = "US$99"
$ws.Range("C11").Value = "'7.4"
$ws.Range("D11").Value = "Good"
$ws.Range("E11").Value = "'5,885"

$ws.Range("A12").Value = "The People - Paris Marais"
$ws.Range("B12").Value = "US$159"
$ws.Range("C12").Value = "'8.3"
$ws.Range("D12").Value = "Very Good"
$ws.Range("E12").Value = "'6,609"

$ws.Range("A13").Value = "Hotel des Carmes by Malone"
$ws.Range("B13").Value = "US$776"
$ws.Range("C13").Value = "'10"
$ws.Range("D13").Value = "Exceptional"
$ws.Range("E13").Value = "'1"

$ws.Range("A14").Value = "Hotel Du Cadran"
$ws.Range("B14").Value = "US$784"
$ws.Range("C14").Value = "'8.5"
$ws.Range("D14").Value = "Very Good"
$ws.Range("E14").Value = "'795"

$ws.Range("A15").Value = "citizenM Paris Champs-Élysées"
$ws.Range("B15").Value = "US$1,167"
$ws.Range("C15").Value = "'8.5"
$ws.Range("D15").Value = "Very Good"
$ws.Range("E15").Value = "'2,787"

$ws.Range("A16").Value = "Motel One Paris-Porte Dorée"
$ws.Range("B16").Value = "US$495"
$ws.Range("C16").Value = "'8.7"
$ws.Range("D16").Value = "Excellent"
$ws.Range("E16").Value = "'6,979"

$ws.Range("A17").Value = "Beau M Paris"
$ws.Range("B17").Value = "US$138"
$ws.Range("C17").Value = "'8.7"
$ws.Range("D17").Value = "Excellent"
$ws.Range("E17").Value = "'1,149"

$ws.Range("A18").Value = "The People Paris Belleville"
$ws.Range("B18").Value = "US$122"
$ws.Range("C18").Value = "'8.2"
$ws.Range("D18").Value = "Very Good"
$ws.Range("E18").Value = "'2,735"

$ws.Range("A19").Value = "Hotel Britannique"
$ws.Range("B19").Value = "US$993"
$ws.Range("C19").Value = "'8.9"
$ws.Range("D19").Value = "Excellent"
$ws.Range("E19").Value = "'1,383"

$ws.Range("A20").Value = "Hotel Relais Bosquet by Malone"
$ws.Range("B20").Value = "US$924"
$ws.Range("C20").Value = "'8.7"
$ws.Range("D20").Value = "Excellent"
$ws.Range("E20").Value = "'1,068"

$ws.Range("A21").Value = "Pullman Paris Montparnasse"
$ws.Range("B21").Value = "US$1,193"
$ws.Range("C21").Value = "'8.4"
$ws.Range("D21").Value = "Very Good"
$ws.Range("E21").Value = "'2,047"

$ws.Range("A22").Value = "Alberte Hôtel"
$ws.Range("B22").Value = "US$989"
$ws.Range("C22").Value = "'9.2"
$ws.Range("D22").Value = "Wonderful"
$ws.Range("E22").Value = "'169"

$ws.Range("A23").Value = "Hôtel Oratio"
$ws.Range("B23").Value = "US$918"
$ws.Range("C23").Value = "'8.8"
$ws.Range("D23").Value = "Excellent"
$ws.Range("E23").Value = "'1,534"

$ws.Range("A24").Value = "Hôtel Henri IV Rive Gauche"
$ws.Range("B24").Value = "US$890"
$ws.Range("C24").Value = "'8.4"
$ws.Range("D24").Value = "Very Good"
$ws.Range("E24").Value = "'1,861"

$ws.Range("A25").Value = "Select Hotel"
$ws.Range("B25").Value = "US$884"
$ws.Range("C25").Value = "'9.0"
$ws.Range("D25").Value = "Wonderful"
$ws.Range("E25").Value = "'2,134"

$ws.Range("A26").Value = "Hotel Ekta Champs Elysées"
$ws.Range("B26").Value = "US$949"
$ws.Range("C26").Value = "'8.8"
$ws.Range("D26").Value = "Excellent"
$ws.Range("E26").Value = "'1,252"

$ws.Range("A27").Value = "Europe Saint Severin-Paris Notre Dame"
$ws.Range("B27").Value = "US$800"
$ws.Range("C27").Value = "'8.4"
$ws.Range("D27").Value = "Very Good"
$ws.Range("E27").Value = "'2,270"
